# 11/05/2021 - Para testing
# Insert 3 new traceability rows at the top of the data (above the existing
# 10/05/2021 rows), shifting everything else down by 3 rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push the existing data (rows 2..24) down to make room for 3 new rows.
$ws.Rows("2:4").Insert()

# The new rows should pick up the plain/default formatting used by the rest
# of the data rows (no inherited borders/number formats from the insert).
$ws.Range("A2:H4").ClearFormats()

# Make sure the date/code columns are written as literal text (matching the
# rest of the sheet, e.g. "00000006" keeps its leading zeros and the date
# stays a string instead of turning into a serial number).
$ws.Range("A2:D4").NumberFormat = "@"

$newRows = @(
  @("11/05/2021", "Egreso", "00000006", "TESTING",        "BLOCK 19X19X39", "UNIDADES", 900,  500),
  @("11/05/2021", "Egreso", "00000006", "TESTING",        "BLOCK 19X19X39", "UNIDADES", 1000, 900),
  @("10/05/2021", "Egreso", "00000005", "MORENO YAMILEM", "BLOCK 19X19X39", "UNIDADES", 1100, 1000)
)

for ($i = 0; $i -lt $newRows.Count; $i++) {
  $r = 2 + $i
  $rowValues = $newRows[$i]
  for ($col = 1; $col -le 8; $col++) {
    $ws.Cells.Item($r, $col).Value = $rowValues[$col - 1]
  }
}

# Drop the temporary text number-format again so the new cells end up with
# the same "no explicit style" look as the surrounding data rows.
$ws.Range("A2:H4").ClearFormats()
